# TMTT0033378 - Verify the functionality of Relationship Treemap added to
# External Contact detail page.
#
# This adds a new "AssociatedEngagements" worksheet (mirroring the existing
# "AffiliatedCompanies" sheet) as the last tab in the workbook, makes it the
# active/selected sheet, and populates it with the four new field labels
# used by the Associated Engagements section.

$wb = $excel.ActiveWorkbook

# Duplicate the existing "AffiliatedCompanies" sheet so the new sheet picks
# up the same look & feel (page margins, default row height, etc.) and lands
# immediately after it as the new last tab.
$source = $wb.Worksheets.Item("AffiliatedCompanies")
$source.Copy($null, $source)

# The copy becomes the last sheet in the workbook - rename it.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "AssociatedEngagements"

# Replace the copied content with the new Associated Engagements field
# labels.
$newSheet.Range("A1").Value = "Engagement:"
$newSheet.Range("A2").Value = "Client Name:"
$newSheet.Range("A3").Value = "Job Type:"
$newSheet.Range("A4").Value = "Role:"

# Match the author's selection/active cell on the new sheet.
$newSheet.Range("F12").Select() | Out-Null

# Make the new sheet the active tab of the workbook.
$newSheet.Activate()
